$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Standardized / Logistic Regression
$ws.Range("C2").Value = "{'C': 100.0, 'solver': 'sag'}"
$ws.Range("D2").Value = 0.7276658484517931

# Row 3 - Standardized / Decision Tree
$ws.Range("D3").Value = 0.729744893757785

# Row 4 - Standardized / Random Forest
$ws.Range("C4").Value = "{'max_depth': 10, 'min_samples_split': 5, 'n_estimators': 50}"
$ws.Range("D4").Value = 0.7346514640225005

# Row 5 - Standardized / Gradient Boosting
$ws.Range("C5").Value = "{'learning_rate': 0.01, 'max_depth': 5, 'n_estimators': 300}"
$ws.Range("D5").Value = 0.734651448893054

# Row 6 - Normalized / Logistic Regression
$ws.Range("D6").Value = 0.7276658484517929

# Row 7 - Normalized / Decision Tree
$ws.Range("C7").Value = "{'criterion': 'entropy', 'max_depth': 8}"
$ws.Range("D7").Value = 0.7297656600037306

# Row 8 - Normalized / Random Forest
$ws.Range("D8").Value = 0.7353375584834156

# Row 9 - Normalized / Gradient Boosting
$ws.Range("C9").Value = "{'learning_rate': 0.01, 'max_depth': 5, 'n_estimators': 300}"
$ws.Range("D9").Value = 0.734672238913844

$wb.Save()
